# Applies the diff described in the task:
#  1) Six "code block" paragraphs get explicit single line spacing
#     (adds w:line="240" to their w:spacing element).
#  2) Three new paragraphs are appended at the end of the body: a blank
#     paragraph, a numbered Heading-3 question paragraph (with a fresh
#     bookmark), and another blank paragraph.
#  3) The section's page margins (top/bottom/left/right) are changed
#     from 1440 twips (1 inch) to ~850 twips (1.5 cm).

$d = $word.ActiveDocument

# --- 1) line spacing on the six code-block paragraphs -----------------
foreach ($i in 4, 5, 6, 7, 8, 9) {
    $p = $d.Paragraphs($i)
    $p.Range.ParagraphFormat.LineSpacingRule = 0   # wdLineSpaceSingle -> w:line="240" w:lineRule="auto"
}

# --- 2) append the three new paragraphs --------------------------------
$countBefore = $d.Paragraphs.Count

$endRange = $d.Paragraphs($countBefore).Range
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

$endRange2 = $d.Paragraphs($d.Paragraphs.Count).Range
$endRange2.Collapse(0)
$endRange2.InsertParagraphAfter()

$endRange3 = $d.Paragraphs($d.Paragraphs.Count).Range
$endRange3.Collapse(0)
$endRange3.InsertParagraphAfter()

# the middle of the three new paragraphs gets the question text + heading
$targetIdx = $countBefore + 2
$newPara = $d.Paragraphs($targetIdx)
$newPara.Range.Text = "escreva um arquivo readme.md para um projeto que inclui os seguintes códigos html e javascript…"

$newPara = $d.Paragraphs($targetIdx)
$newPara.Style = "Heading 3"

# re-use the existing numbered-list definition (numId 1) already used by
# the other question headings in the document
$srcTemplate = $d.Paragraphs(17).Range.ListFormat.ListTemplate
$newPara.Range.ListFormat.ApplyListTemplateWithLevel($srcTemplate, $true, 1, $false, 1)

$newPara.Range.ParagraphFormat.LeftIndent = 36      # 720 twips
$newPara.Range.ParagraphFormat.FirstLineIndent = -18 # -360 twips (hanging)

$bmRange = $d.Range($newPara.Range.Start, $newPara.Range.Start)
$d.Bookmarks.Add("_15fbj9s3uuq4", $bmRange)

# --- 3) page margins: 1440 twips (1in) -> ~850 twips (1.5cm) ----------
$ps = $d.PageSetup
$ps.TopMargin = 850.3937007874016 / 20
$ps.BottomMargin = 850.3937007874016 / 20
$ps.LeftMargin = 850.3937007874016 / 20
$ps.RightMargin = 850.3937007874016 / 20

Write-Output "edit complete: paragraphs=$($d.Paragraphs.Count)"
